$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells H1:J1 ("Unnamed: 7/8/9" - artifact of a pandas re-export) ---
$ws.Range("H1").Value = "Unnamed: 7"
$ws.Range("I1").Value = "Unnamed: 8"
$ws.Range("J1").Value = "Unnamed: 9"
# Match the header formatting (bold, bordered, centered) already used by A1:G1
$ws.Range("G1").Copy()
$ws.Range("H1:J1").PasteSpecial(-4122)

# --- Turn the three summary formulas (J2:J4) into their last computed values ---
$ws.Range("J2").Value = 13
$ws.Range("J3").Value = 22
# Drop the percent number-format on J4 before writing the plain decimal value
$ws.Range("J4").Style = "Normal"
$ws.Range("J4").Value = 0.59090909090909094

# --- Fill in the "Beat Vegas?" results for the Jan 7 games that were previously blank ---
$ws.Range("G18").Value = "No"
$ws.Range("G19").Value = "No"
$ws.Range("G20").Value = "No"
$ws.Range("G21").Value = "No"
$ws.Range("G22").Value = "Yes"

# --- Append the Jan 8 games (rows 23-32), pasting the existing date format first ---
$ws.Range("A2").Copy()
$ws.Range("A23:A32").PasteSpecial(-4122)

$games = @(
  @(44204, "DET", "PHO", 8,    9.4,   -1.4),
  @(44204, "NOP", "CHO", -7,   -1.6,  -5.4),
  @(44204, "BOS", "WAS", -4.5, 3.4,   -7.9),
  @(44204, "NYK", "OKC", -2.5, 2.5,   -5),
  @(44204, "MIL", "UTA", -5.5, -20.3, 14.8),
  @(44204, "HOU", "ORL", -6.5, -5.2,  -1.3),
  @(44204, "MEM", "BRK", -1,   16.7,  -17.7),
  @(44204, "SAC", "TOR", 5,    -6.4,  11.4),
  @(44204, "GSW", "LAC", 6.5,  -3.9,  10.4),
  @(44204, "LAL", "CHI", -9.5, 1,     -10.5)
)

$r = 23
foreach ($g in $games) {
  $ws.Cells.Item($r, 1).Value = $g[0]
  $ws.Cells.Item($r, 2).Value = $g[1]
  $ws.Cells.Item($r, 3).Value = $g[2]
  $ws.Cells.Item($r, 4).Value = $g[3]
  $ws.Cells.Item($r, 5).Value = $g[4]
  $ws.Cells.Item($r, 6).Value = $g[5]
  $r = $r + 1
}
